$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report row is inserted at row 187 (pushing the existing
# historical rows 187-314 down to 188-315). This is equivalent to inserting
# a new row at the top of the data block and keeping everything else intact.
$ws.Rows(187).Insert()

# Fill in the newly inserted row 187 with this week's values. All of the
# "constant" columns (A,B,C,E,F,G,H,I,R) share the same values throughout
# the whole data set, so we just repeat them here as well.
$ws.Cells.Item(187, 1).Value = 5
$ws.Cells.Item(187, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(187, 3).Value = "Maule"
$ws.Cells.Item(187, 4).Value = 44824
$ws.Cells.Item(187, 5).Value = 7
$ws.Cells.Item(187, 6).Value = 100112009
$ws.Cells.Item(187, 7).Value = "Acelga"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 400
$ws.Cells.Item(187, 11).Value = 2500
$ws.Cells.Item(187, 12).Value = 2500
$ws.Cells.Item(187, 13).Value = 2500
$ws.Cells.Item(187, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(187, 15).Value = "Región del Maule"
$ws.Cells.Item(187, 16).Value = 625
$ws.Cells.Item(187, 17).Value = 4
$ws.Cells.Item(187, 18).Value = "Hortaliza"
